$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152 - this shifts existing rows 152..251 down to 153..252
# and grows the sheet dimension to A1:R252 automatically.
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with the new weekly data point.
$ws.Range("A152").Value = 7
$ws.Range("B152").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C152").Value = "Ñuble"
$ws.Range("D152").Value = 44596
$ws.Range("E152").Value = 16
$ws.Range("F152").Value = 100114013
$ws.Range("G152").Value = "Zanahoria"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Segunda"
$ws.Range("J152").Value = 100
$ws.Range("K152").Value = 5000
$ws.Range("L152").Value = 5500
$ws.Range("M152").Value = 5250
$ws.Range("N152").Value = "`$/saco 20 kilos"
$ws.Range("O152").Value = "Provincia de Diguillín"
$ws.Range("P152").Value = 262
$ws.Range("Q152").Value = 20
$ws.Range("R152").Value = "Hortaliza"

# Match the date-cell number format used throughout column D.
$ws.Range("D152").NumberFormat = $ws.Range("D153").NumberFormat
